$wb = $excel.ActiveWorkbook

# --- Carrier sheet: just move the selected/active cell (cursor moved to A6) ---
$wsCarrier = $wb.Worksheets.Item("Carrier")
$wsCarrier.Range("A6").Select()

# --- Generator sheet: add a p_max_pu column and a new "diesel" generator row ---
$ws = $wb.Worksheets.Item("Generator")

# New header for column F
$ws.Range("F1").Value = "p_max_pu"

# New data row (row 3) describing the diesel generator
$ws.Range("A3").Value = "diesel"
$ws.Range("B3").Value = "AC"
$ws.Range("C3").Value = "bus 0"

# D2 already holds the text "True" (not boolean) - copy it down so D3 keeps
# the same shared-string text type instead of becoming a native boolean.
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4163)

$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1

# Column F should auto-size like the other "bestFit" columns
$ws.Columns.Item(6).AutoFit()

# Move the active cell on the Generator sheet to A4
$ws.Range("A4").Select()
